$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 24240578
$ws.Range("I19").Value = 18783524
$ws.Range("K19").Value = 18783524
$ws.Range("M19").Value = -18783349

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3334
$ws.Range("I32").Value = 2250
$ws.Range("J32").Value = 3876
$ws.Range("K32").Value = 2250
$ws.Range("L32").Value = 3876
$ws.Range("M32").Value = -1924
$ws.Range("N32").Value = -4528

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 638.1
$ws.Range("J43").Value = 738
$ws.Range("L43").Value = 738
$ws.Range("N43").Value = -876

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3966.608
$ws.Range("J64").Value = 4111.5
$ws.Range("L64").Value = 4111.5
$ws.Range("N64").Value = -4607.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3966.608
$ws.Range("J67").Value = 4111.5
$ws.Range("L67").Value = 4111.5
$ws.Range("N67").Value = -5827.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 898.8283
$ws.Range("J129").Value = 898.8283
$ws.Range("L129").Value = 2696.4849
$ws.Range("N129").Value = -12696.4849

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21744060
$ws.Range("I32").Value = 28574090
$ws.Range("K32").Value = 28574090
$ws.Range("M32").Value = -28573803

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6168.75
$ws.Range("I61").Value = 2234
$ws.Range("K61").Value = 2234
$ws.Range("M61").Value = -2022

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 29499.5
$ws.Range("J76").Value = 29499.5
$ws.Range("L76").Value = 29499.5
$ws.Range("N76").Value = -30175.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 29499.5
$ws.Range("J79").Value = 29499.5
$ws.Range("L79").Value = 29499.5
$ws.Range("N79").Value = -31839.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2144.6667
$ws.Range("J88").Value = 2233.3333
$ws.Range("L88").Value = 2233.3333
$ws.Range("N88").Value = -3045.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2144.6667
$ws.Range("J91").Value = 2233.3333
$ws.Range("L91").Value = 2233.3333
$ws.Range("N91").Value = -5041.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1524.6666
$ws.Range("I122").Value = 1340.25
$ws.Range("K122").Value = 4020.75
$ws.Range("M122").Value = -1570.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6168.75
$ws.Range("I136").Value = 2234
$ws.Range("K136").Value = 6702
$ws.Range("M136").Value = -4152

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 4825
$ws.Range("J58").Value = 4825
$ws.Range("L58").Value = 4825
$ws.Range("N58").Value = -5413

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 37163.83
$ws.Range("J59").Value = 37163.83
$ws.Range("L59").Value = 37163.83
$ws.Range("N59").Value = -38857.83

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 40780
$ws.Range("J60").Value = 40780
$ws.Range("L60").Value = 40780
$ws.Range("N60").Value = -41978

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 16279.429
$ws.Range("I82").Value = 5789
$ws.Range("J82").Value = 30266.666
$ws.Range("K82").Value = 5789
$ws.Range("L82").Value = 30266.666
$ws.Range("M82").Value = -5406
$ws.Range("N82").Value = -31032.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 16279.429
$ws.Range("I85").Value = 5789
$ws.Range("J85").Value = 30266.666
$ws.Range("K85").Value = 5789
$ws.Range("L85").Value = 30266.666
$ws.Range("M85").Value = -4463
$ws.Range("N85").Value = -32918.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2625
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2625
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2625
$ws.Range("N86").Value = -4871
$ws.Range("M86").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2625
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2625
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 13125
$ws.Range("N89").Value = -24357
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2430.8572
$ws.Range("I107").Value = 2581.9
$ws.Range("J107").Value = 2053.25
$ws.Range("K107").Value = 2581.9
$ws.Range("L107").Value = 2053.25
$ws.Range("M107").Value = -661.9000000000001
$ws.Range("N107").Value = -5893.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2021.2
$ws.Range("I134").Value = 2021.2
$ws.Range("K134").Value = 6063.6
$ws.Range("M134").Value = -3528.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2975.25
$ws.Range("I58").Value = 2717
$ws.Range("J58").Value = 3750
$ws.Range("K58").Value = 2717
$ws.Range("L58").Value = 3750
$ws.Range("M58").Value = -2514
$ws.Range("N58").Value = -4156

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1671.375
$ws.Range("I105").Value = 1671.375
$ws.Range("K105").Value = 1671.375
$ws.Range("M105").Value = 75.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 18900
$ws.Range("J109").Value = 18900
$ws.Range("L109").Value = 18900
$ws.Range("N109").Value = -20980

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2975.25
$ws.Range("I136").Value = 2717
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 8151
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -5601
$ws.Range("N136").Value = -16350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 744.0833
$ws.Range("I5").Value = 616.8
$ws.Range("J5").Value = 835
$ws.Range("K5").Value = 1850.4
$ws.Range("L5").Value = 2505
$ws.Range("M5").Value = -1738.4
$ws.Range("N5").Value = -2729

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 2764.1428
$ws.Range("I115").Value = 1639.6
$ws.Range("J115").Value = 3388.889
$ws.Range("K115").Value = 4918.799999999999
$ws.Range("L115").Value = 10166.667
$ws.Range("M115").Value = -3743.799999999999
$ws.Range("N115").Value = -12516.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 8300.643
$ws.Range("I122").Value = 11015.95
$ws.Range("K122").Value = 99143.55
$ws.Range("M122").Value = -96693.55

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 5905.2085
$ws.Range("I129").Value = 1903.8
$ws.Range("J129").Value = 8763.357
$ws.Range("K129").Value = 5711.4
$ws.Range("L129").Value = 26290.071
$ws.Range("M129").Value = -711.3999999999996
$ws.Range("N129").Value = -36290.071

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 878.9259
$ws.Range("J131").Value = 1036.6471
$ws.Range("L131").Value = 3109.9413
$ws.Range("N131").Value = -13189.9413

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 744.0833
$ws.Range("I135").Value = 616.8
$ws.Range("J135").Value = 835
$ws.Range("K135").Value = 5551.2
$ws.Range("L135").Value = 7515
$ws.Range("M135").Value = -3016.2
$ws.Range("N135").Value = -12585

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4473.1904
$ws.Range("I132").Value = 4218.125
$ws.Range("K132").Value = 12654.375
$ws.Range("M132").Value = -10124.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000860.2
$ws.Range("I22").Value = 1429128.6
$ws.Range("J22").Value = 1567.3334
$ws.Range("K22").Value = 1429128.6
$ws.Range("L22").Value = 1567.3334
$ws.Range("M22").Value = -1428833.6
$ws.Range("N22").Value = -2157.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1000860.2
$ws.Range("I27").Value = 1429128.6
$ws.Range("J27").Value = 1567.3334
$ws.Range("K27").Value = 1429128.6
$ws.Range("L27").Value = 1567.3334
$ws.Range("M27").Value = -1429021.6
$ws.Range("N27").Value = -1781.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3399.75
$ws.Range("I136").Value = 1882
$ws.Range("J136").Value = 4917.5
$ws.Range("K136").Value = 5646
$ws.Range("L136").Value = 14752.5
$ws.Range("M136").Value = -3096
$ws.Range("N136").Value = -19852.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 35904
$ws.Range("J137").Value = 35904
$ws.Range("L137").Value = 35904
$ws.Range("N137").Value = -46104

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7632
$ws.Range("J41").Value = 7632
$ws.Range("L41").Value = 7632
$ws.Range("N41").Value = -8412

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 17729
$ws.Range("J45").Value = 17638.666
$ws.Range("L45").Value = 17638.666
$ws.Range("N45").Value = -18620.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2094.182
$ws.Range("I122").Value = 1774.4
$ws.Range("J122").Value = 2586.1538
$ws.Range("K122").Value = 5323.200000000001
$ws.Range("L122").Value = 7758.4614
$ws.Range("M122").Value = -2873.200000000001
$ws.Range("N122").Value = -12658.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 29850
$ws.Range("J133").Value = 29850
$ws.Range("L133").Value = 29850
$ws.Range("N133").Value = -39970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 15870
$ws.Range("I136").Value = 20160
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 60480
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -57930
$ws.Range("N136").Value = -14100
